# Add a new "Reg Proc" column to the requirements table on the "Details"
# sheet, and populate the clarification notes that go with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Details")
$lo = $ws.ListObjects.Item("Table2")

# Add the new table column (this extends the table range and autofilter,
# and shifts the used range out to column T).
$newCol = $lo.ListColumns.Add()

# Match the header formatting of the other header cells in row 2 (the new
# cell otherwise inherits the row's default/custom format instead of the
# explicit per-cell header style).
$ws.Range("S2").Copy()
$ws.Range("T2").PasteSpecial(-4122)
$newCol.Range.Cells(1, 1).Value = "Reg Proc"

# Populate the new column's notes.
$ws.Range("T4").Value = "When UIN IS needed to be generated`n1.the Acknowledgment from Print queue- what needs to be done`nTime period `n2. If there is a print failure- no need to handle from MOSIP`nUser Story ?"
$ws.Range("T5").Value = "No Mapping of such kind from Reg Processor`nID Repo- Might not be there in ID Repo as well"
$ws.Range("T6").Value = "ID Repo- need to know "
$ws.Range("T7").Value = "there shud be a label as Res_Service`nReg Client packet needs to be understood`nService from Reg proc needs to be developed"
$ws.Range("T8").Value = "Under processing`nProcessed`n"
$ws.Range("T9").Value = "Under processing`nProcessed"
$ws.Range("T10").Value = "E-UIN Generation"

# Existing "Comments" note in row 8 gets an addendum.
$ws.Range("S8").Value = "Reg proc`nArchival policy"

# Wrap text on the new notes to match the rest of the sheet's commentary
# cells, and column width to fit.
$ws.Range("T4").WrapText = $true
$ws.Range("T5").WrapText = $true
$ws.Range("T7").WrapText = $true
$ws.Range("T8").WrapText = $true
$ws.Range("T9").WrapText = $true
$ws.Range("S8").WrapText = $true
$ws.Columns.Item(20).ColumnWidth = 32.08984375

# Keep the view roughly where the author left it.
$ws.Range("T4").Select()
